# "updated assumptions and constraints"
#
# - A27 ("Steel tile with") -> "Steel Tile"
# - A30 ("screed")          -> "Screed"
# - Leave the active selection on H22 (matches the saved selection in the diff)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("assumptions")

$ws.Range("A27").Value = "Steel Tile"
$ws.Range("A30").Value = "Screed"

$ws.Range("H22").Select()
